# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Thu Jun 15 20:56:59 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.353.79"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3
$ws.Range("D3").Value = "1.660.17"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.79%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.43"
$ws.Range("E5").Value = "  -0.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.65%  "

# Row 7
$ws.Range("E7").Value = "  +1.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2607"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06137"
$ws.Range("E9").Value = "  +2.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07057"
$ws.Range("E10").Value = "  -0.85%  "

# Row 11
$ws.Range("D11").Value = "1.659.66"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.68"
$ws.Range("E12").Value = "  +2.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5929"
$ws.Range("E13").Value = "  -4.50%  "

# Row 14
$ws.Range("E14").Value = "  -5.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.34"
$ws.Range("E15").Value = "  +1.92%  "

# Row 16
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").Value = "  -0.82%  "

# Row 18
$ws.Range("D18").Value = "25.334.34"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006765"
$ws.Range("E19").Value = "  +3.03%  "

# Row 20
$ws.Range("E20").Value = "  -0.08%  "

# Row 21
$ws.Range("D21").Value = "1.872.51"
$ws.Range("E21").Value = "  -1.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.451"
$ws.Range("E22").Value = "  +0.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.630"
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.337"
$ws.Range("E24").Value = "  +1.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "133.88"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.99"
$ws.Range("E26").Value = "  +1.94%  "

# Row 27
$ws.Range("E27").Value = "  +2.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "103.89"
$ws.Range("E28").Value = "  +1.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.682"
$ws.Range("E29").Value = "  -1.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.955"
$ws.Range("E30").Value = "  +3.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.630"
$ws.Range("E31").Value = "  +2.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07640"
$ws.Range("E32").Value = "  -3.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04336"
$ws.Range("E33").Value = "  -5.88%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9997"
$ws.Range("E34").Value = "  -0.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.605"
$ws.Range("E35").Value = "  -1.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6161"
$ws.Range("E36").Value = "  +6.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9450"
$ws.Range("E37").Value = "  +0.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.601"
$ws.Range("E38").Value = "  -1.08%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8566"
$ws.Range("E39").Value = "  +2.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.84%  "

# Row 41
$ws.Range("E41").Value = "  -2.98%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.830"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.62"
$ws.Range("E43").Value = "  -1.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3770"
$ws.Range("E44").Value = "  +1.86%  "

# Row 45
$ws.Range("E45").Value = "  -4.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1114"
$ws.Range("E46").Value = "  -1.85%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.181"
$ws.Range("E47").Value = "  +2.37%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05247"
$ws.Range("E48").Value = "  +1.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.42"
$ws.Range("E49").Value = "  -0.78%  "

# Row 50
$ws.Range("E50").Value = "  -0.65%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.324"
$ws.Range("E51").Value = "  -0.48%  "
